$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2486.0557
$ws.Range("I18").Value = 2783.2666
$ws.Range("K18").Value = 2783.2666
$ws.Range("M18").Value = -2499.2666

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 808.5909
$ws.Range("I28").Value = 623.5
$ws.Range("J28").Value = 1302.1666
$ws.Range("K28").Value = 623.5
$ws.Range("L28").Value = 1302.1666
$ws.Range("M28").Value = -138.5
$ws.Range("N28").Value = -2272.1666

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1109.6666
$ws.Range("I98").Value = 1136.2222
$ws.Range("J98").Value = 1030
$ws.Range("K98").Value = 1136.2222
$ws.Range("L98").Value = 1030
$ws.Range("M98").Value = 361.7778000000001
$ws.Range("N98").Value = -4026

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 15272.857
$ws.Range("I113").Value = 9466.666999999999
$ws.Range("J113").Value = 16856.363
$ws.Range("K113").Value = 9466.666999999999
$ws.Range("L113").Value = 16856.363
$ws.Range("M113").Value = -6212.666999999999
$ws.Range("N113").Value = -23364.363

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7619.579
$ws.Range("I116").Value = 4964.6665
$ws.Range("K116").Value = 4964.6665
$ws.Range("M116").Value = -1522.6665

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1109.6666
$ws.Range("I122").Value = 1136.2222
$ws.Range("J122").Value = 1030
$ws.Range("K122").Value = 3408.6666
$ws.Range("L122").Value = 3090
$ws.Range("M122").Value = -958.6665999999996
$ws.Range("N122").Value = -7990

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 203821.3
$ws.Range("I129").Value = 288044.72
$ws.Range("J129").Value = 7300
$ws.Range("K129").Value = 864134.1599999999
$ws.Range("L129").Value = 21900
$ws.Range("M129").Value = -859134.1599999999
$ws.Range("N129").Value = -31900

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 26733.682
$ws.Range("I132").Value = 1790.25
$ws.Range("K132").Value = 5370.75
$ws.Range("M132").Value = -2840.75

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3001.1777
$ws.Range("I138").Value = 1383.0769
$ws.Range("J138").Value = 3658.5312
$ws.Range("K138").Value = 4149.2307
$ws.Range("L138").Value = 10975.5936
$ws.Range("M138").Value = 990.7692999999999
$ws.Range("N138").Value = -21255.5936

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7047.6665
$ws.Range("I141").Value = 6131.923
$ws.Range("K141").Value = 18395.769
$ws.Range("M141").Value = -13215.769

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7717
$ws.Range("I45").Value = 4498
$ws.Range("K45").Value = 4498
$ws.Range("M45").Value = -4121

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17147.111
$ws.Range("I61").Value = 2514.8
$ws.Range("K61").Value = 2514.8
$ws.Range("M61").Value = -2302.8

# ARM row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 17147.111
$ws.Range("I136").Value = 2514.8
$ws.Range("K136").Value = 7544.400000000001
$ws.Range("M136").Value = -4994.400000000001

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2418164.5
$ws.Range("I94").Value = 1983.4736
$ws.Range("K94").Value = 1983.4736
$ws.Range("M94").Value = -1532.4736

# CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 468.08334
$ws.Range("I10").Value = 367.66666
$ws.Range("J10").Value = 769.3333
$ws.Range("K10").Value = 367.66666
$ws.Range("L10").Value = 769.3333
$ws.Range("M10").Value = -228.66666
$ws.Range("N10").Value = -1047.3333

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 58215.625
$ws.Range("J86").Value = 28645.2
$ws.Range("L86").Value = 28645.2
$ws.Range("N86").Value = -30891.2

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 58215.625
$ws.Range("J89").Value = 28645.2
$ws.Range("L89").Value = 143226
$ws.Range("N89").Value = -154458

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 9972.166999999999
$ws.Range("I107").Value = 815.6667
$ws.Range("J107").Value = 37441.668
$ws.Range("K107").Value = 815.6667
$ws.Range("L107").Value = 37441.668
$ws.Range("M107").Value = 1104.3333
$ws.Range("N107").Value = -41281.668

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6200.222
$ws.Range("I134").Value = 6349.375
$ws.Range("K134").Value = 19048.125
$ws.Range("M134").Value = -16513.125

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1478.6471
$ws.Range("I131").Value = 982.5
$ws.Range("K131").Value = 2947.5
$ws.Range("M131").Value = 2092.5

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1914.2858
$ws.Range("I132").Value = 1466.6666
$ws.Range("K132").Value = 13199.9994
$ws.Range("M132").Value = -10669.9994

# GSM row 7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 6700000
$ws.Range("J7").Value = 6700000
$ws.Range("L7").Value = 6700000
$ws.Range("N7").Value = -6700224

# GSM row 8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 6700000
$ws.Range("J8").Value = 6700000
$ws.Range("L8").Value = 6700000
$ws.Range("N8").Value = -6700278

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 891333.3
$ws.Range("I11").Value = 50000
$ws.Range("J11").Value = 1312000
$ws.Range("K11").Value = 50000
$ws.Range("L11").Value = 1312000
$ws.Range("M11").Value = -49861
$ws.Range("N11").Value = -1312278

# GSM row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# GSM row 19
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 2199.75
$ws.Range("I19").Value = 1500
$ws.Range("J19").Value = 2433
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 2433
$ws.Range("M19").Value = -1212
$ws.Range("N19").Value = -3009

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6994023
$ws.Range("I22").Value = 12987701
$ws.Range("K22").Value = 12987701
$ws.Range("M22").Value = -12987406

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 6994023
$ws.Range("I27").Value = 12987701
$ws.Range("K27").Value = 12987701
$ws.Range("M27").Value = -12987594

# LTW row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# LTW row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3770.7273
$ws.Range("I122").Value = 3770.7273
$ws.Range("K122").Value = 11312.1819
$ws.Range("M122").Value = -8862.1819

# WVR row 7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4583.222
$ws.Range("I132").Value = 2984.963
$ws.Range("K132").Value = 8954.889000000001
$ws.Range("M132").Value = -6424.889000000001
